# add rural communities interactions back in, fix de_dg files (not run yet)
#
# The "Legislature" row (row 12, column A) is removed from the sheet.
# Deleting the entire row shifts all subsequent rows up by one and
# automatically removes the now-unused "Legislature" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Delete()
